# Updates to C14 table: add a comma separator between the context type
# (e.g. "burial"/"midden") and the feature code (e.g. "M020"/"H044")
# in the "Context" column (column I) of Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("I10").Value = "burial, M020"
$ws.Range("I11").Value = "midden, H044"
$ws.Range("I12").Value = "burial, M009"
$ws.Range("I16").Value = "burial, M039"
$ws.Range("I19").Value = "burial, M066"
$ws.Range("I21").Value = "midden, H193"
$ws.Range("I23").Value = "midden, H026"
$ws.Range("I24").Value = "burial, M095"

# Update the saved view state: selection moves to F1 and the view is no
# longer scrolled down (topLeftCell reset back to default).
$ws.Range("F1").Select() | Out-Null
